$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Element" values between row 2 and row 3 (column B):
# B2 currently "2-sheets (Collated/Nested)  2p" -> should become "1-sheet (Collated/Nested)  2p"
# B3 currently "1-sheet (Collated/Nested)  2p" -> should become "2-sheets (Collated/Nested)  2p"
$ws.Range("B2").Value = "1-sheet (Collated/Nested)  2p"
$ws.Range("B3").Value = "2-sheets (Collated/Nested)  2p"
